$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1675.4166
$ws.Range("I100").Value = 1350.8334
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1350.8334
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -809.8334
$ws.Range("N100").Value = -3082
$ws.Range("H109").Value = 40649.5
$ws.Range("J109").Value = 40649.5
$ws.Range("L109").Value = 40649.5
$ws.Range("N109").Value = -43423.5
$ws.Range("H117").Value = 48742
$ws.Range("J117").Value = 48742
$ws.Range("L117").Value = 48742
$ws.Range("N117").Value = -57920
$ws.Range("H124").Value = 41878.715
$ws.Range("J124").Value = 46358.5
$ws.Range("L124").Value = 46358.5
$ws.Range("N124").Value = -56178.5
$ws.Range("H128").Value = 43172
$ws.Range("J128").Value = 43172
$ws.Range("L128").Value = 43172
$ws.Range("N128").Value = -53132

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1990.3667
$ws.Range("I74").Value = 1995.2727
$ws.Range("J74").Value = 1976.875
$ws.Range("K74").Value = 1995.2727
$ws.Range("L74").Value = 1976.875
$ws.Range("M74").Value = -1121.2727
$ws.Range("N74").Value = -3724.875
$ws.Range("H77").Value = 1990.3667
$ws.Range("I77").Value = 1995.2727
$ws.Range("J77").Value = 1976.875
$ws.Range("K77").Value = 9976.363499999999
$ws.Range("L77").Value = 9884.375
$ws.Range("M77").Value = -5608.363499999999
$ws.Range("N77").Value = -18620.375
$ws.Range("H117").Value = 48415.332
$ws.Range("J117").Value = 48415.332
$ws.Range("L117").Value = 48415.332
$ws.Range("N117").Value = -57593.332
$ws.Range("H118").Value = 49206
$ws.Range("J118").Value = 49206
$ws.Range("L118").Value = 49206
$ws.Range("N118").Value = -52520
$ws.Range("H125").Value = 46450.832
$ws.Range("J125").Value = 46450.832
$ws.Range("L125").Value = 46450.832
$ws.Range("N125").Value = -56290.832
$ws.Range("H131").Value = 45488.168
$ws.Range("J131").Value = 45488.168
$ws.Range("L131").Value = 45488.168
$ws.Range("N131").Value = -55568.168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 47676
$ws.Range("J108").Value = 47676
$ws.Range("L108").Value = 47676
$ws.Range("N108").Value = -55356
$ws.Range("H110").Value = 48694
$ws.Range("J110").Value = 48694
$ws.Range("L110").Value = 48694
$ws.Range("N110").Value = -56874
$ws.Range("H117").Value = 49742
$ws.Range("J117").Value = 49742
$ws.Range("L117").Value = 49742
$ws.Range("N117").Value = -58920
$ws.Range("H124").Value = 52996
$ws.Range("J124").Value = 52996
$ws.Range("L124").Value = 52996
$ws.Range("N124").Value = -62816
$ws.Range("H125").Value = 50780
$ws.Range("J125").Value = 50780
$ws.Range("L125").Value = 50780
$ws.Range("N125").Value = -60620
$ws.Range("H126").Value = 50780
$ws.Range("J126").Value = 50780
$ws.Range("L126").Value = 50780
$ws.Range("N126").Value = -60660
$ws.Range("H130").Value = 49383.5
$ws.Range("J130").Value = 49383.5
$ws.Range("L130").Value = 49383.5
$ws.Range("N130").Value = -59423.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 37858.5
$ws.Range("J110").Value = 37858.5
$ws.Range("L110").Value = 37858.5
$ws.Range("N110").Value = -46038.5
$ws.Range("H111").Value = 48694
$ws.Range("J111").Value = 48694
$ws.Range("L111").Value = 48694
$ws.Range("N111").Value = -56874
$ws.Range("H112").Value = 31956.666
$ws.Range("J112").Value = 31956.666
$ws.Range("L112").Value = 31956.666
$ws.Range("N112").Value = -34910.666
$ws.Range("H116").Value = 49368.5
$ws.Range("J116").Value = 49368.5
$ws.Range("L116").Value = 49368.5
$ws.Range("N116").Value = -58546.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 47659.332
$ws.Range("J110").Value = 47659.332
$ws.Range("L110").Value = 47659.332
$ws.Range("N110").Value = -55839.332
$ws.Range("H116").Value = 39000
$ws.Range("J116").Value = 39000
$ws.Range("L116").Value = 39000
$ws.Range("N116").Value = -48178
$ws.Range("H119").Value = 34380.5
$ws.Range("J119").Value = 34380.5
$ws.Range("L119").Value = 34380.5
$ws.Range("N119").Value = -44056.5
$ws.Range("H130").Value = 45391.332
$ws.Range("J130").Value = 45391.332
$ws.Range("L130").Value = 45391.332
$ws.Range("N130").Value = -55431.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
$ws.Range("H114").Value = 27592
$ws.Range("J114").Value = 27592
$ws.Range("L114").Value = 27592
$ws.Range("N114").Value = -36270
$ws.Range("H118").Value = 40702.5
$ws.Range("J118").Value = 40702.5
$ws.Range("L118").Value = 40702.5
$ws.Range("N118").Value = -44016.5
$ws.Range("H124").Value = 32819.75
$ws.Range("J124").Value = 32819.75
$ws.Range("L124").Value = 32819.75
$ws.Range("N124").Value = -42639.75
$ws.Range("H125").Value = 46426.75
$ws.Range("J125").Value = 46426.75
$ws.Range("L125").Value = 46426.75
$ws.Range("N125").Value = -56266.75
$ws.Range("H127").Value = 46530.11
$ws.Range("J127").Value = 46530.11
$ws.Range("L127").Value = 46530.11
$ws.Range("N127").Value = -56450.11
$ws.Range("H128").Value = 47421
$ws.Range("J128").Value = 47421
$ws.Range("L128").Value = 47421
$ws.Range("N128").Value = -57381

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 44618
$ws.Range("J108").Value = 44618
$ws.Range("L108").Value = 44618
$ws.Range("N108").Value = -52298
$ws.Range("H110").Value = 48644
$ws.Range("J110").Value = 48644
$ws.Range("L110").Value = 48644
$ws.Range("N110").Value = -56824
$ws.Range("H117").Value = 44068
$ws.Range("J117").Value = 44068
$ws.Range("L117").Value = 44068
$ws.Range("N117").Value = -53246
$ws.Range("H120").Value = 46420
$ws.Range("J120").Value = 46420
$ws.Range("L120").Value = 46420
$ws.Range("N120").Value = -56096
$ws.Range("H131").Value = 50715
$ws.Range("J131").Value = 50715
$ws.Range("L131").Value = 50715
$ws.Range("N131").Value = -60795
